# Apply the odds updates described in the diff for 2026-02-16 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 1.66

$ws.Range("G3").Value = 6.4
$ws.Range("H3").Value = 1.76

$ws.Range("I4").Value = 7

$ws.Range("F5").Value = 1.44
$ws.Range("G5").Value = 1.52
$ws.Range("H5").Value = 9.199999999999999
$ws.Range("I5").Value = 12
$ws.Range("K5").Value = 4.8
$ws.Range("Q5").Value = 2.22

$ws.Range("H9").Value = 3.65
$ws.Range("J9").Value = 4
$ws.Range("P9").Value = 2.3

$ws.Range("H10").Value = 1.96
$ws.Range("K10").Value = 6.2

$ws.Range("J11").Value = 2.98
$ws.Range("K11").Value = 5.2

$ws.Range("F13").Value = 1.23
$ws.Range("G13").Value = 1.4
$ws.Range("H13").Value = 3.6
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 5.6
$ws.Range("K13").Value = 980
$ws.Range("N13").Value = 1.02
$ws.Range("O13").Value = 1.24
$ws.Range("Q13").Value = 1.24

$ws.Range("I15").Value = 1.6
$ws.Range("J15").Value = 4.5

$ws.Range("J16").Value = 3.8

$ws.Range("H21").Value = 2.54
$ws.Range("I21").Value = 2.76

$ws.Range("S24").Value = 1.9

$wb.Save()
